# Append weekly pricing rows (240-246) to the Girasol daily pricing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45950, 25966.02, 4700, 4225, 26441.02, 601.01,  0,  0, 601.01,  27042.03),
    @(45951, 44408.7,   330, 1059, 43679.7,   349.7,  50,  0, 399.7,   44079.39999999999),
    @(45952, 72082.62,  110,  600, 71592.62,  707.14,  0,  0, 707.14,  72299.75999999999),
    @(45953, 56410.02,  620,  700, 56330.02, 1236,     0,  0, 1236,    57566.02),
    @(45954, 46762.53,   60,  410, 46412.53, 2406.62,  0,  0, 2406.62, 48819.15),
    @(45955, 30,          0,    0, 30,          0,     0,  0, 0,       30),
    @(45957, 950,         0,    0, 950,         0,     0,  0, 0,       950)
)

$startRow = 240
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = "GIRASOL"

    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
}
